$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the values for the new "House Robber II" row first.
$ws.Range("A20").Value = "House Robber II - Houses in Circle"
$ws.Range("B20").Value = "Return maximum ammount without robbing adjacent houses"
$ws.Range("C20").Value = "Split input array into two arrays skipping first and last houses. Pass these two arrays to the house robber function and return max value"
$ws.Range("D20").Value = "https://leetcode.com/problems/house-robber-ii/"

# Register the hyperlink for the new link cell.
$ws.Hyperlinks.Add($ws.Range("D20"), "https://leetcode.com/problems/house-robber-ii/")

# Copy formatting from the previous data row (row 19) down into the new row 20
# so the new row picks up the same cell styles (column A uses the "Neutral"
# style, B/C use default, D uses the "Hyperlink" style) without leaving extra
# unused style entries attached to the new cells.
$ws.Range("A19:D19").Copy()
$ws.Range("A20:D20").PasteSpecial(-4122)

# Update the active selection as recorded in the saved view state.
$ws.Range("C13").Select()
